# The "historia" paragraph that opens with "Você é Trevor Smith..." is
# rewritten: the short "O ano é 2066, dois anos após o falecimento..."
# sentence becomes a fuller recap of Kat's death and the paragraph's old
# closing ("As tardes quentes ... mais escuras.") is dropped entirely and
# replaced with new text about Katelyn "Kat" Parker and Trevor's quest for
# peace. The trailing _GoBack bookmark is kept, now at the very end of the
# paragraph.

$d = $word.ActiveDocument

# Locate the target paragraph reliably by its distinctive content instead
# of a hard-coded index.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Você é*" -and $t -like "*O ano*2066*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the 'Você é Trevor Smith / O ano é 2066' paragraph"
}

$r = $d.Paragraphs.Item($targetIndex).Range

# Common run formatting used throughout this paragraph.
$rPr = '<w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'

# NOTE: this PowerShell host does not honour default parameter values
# (an omitted argument binds to $null/empty instead of the declared
# default), so every call below passes $true/$false explicitly.
function Run([string]$text, [bool]$preserve) {
    if ($preserve) {
        return '<w:r>' + $rPr + '<w:t xml:space="preserve">' + $text + '</w:t></w:r>'
    } else {
        return '<w:r>' + $rPr + '<w:t>' + $text + '</w:t></w:r>'
    }
}

$spellStart = '<w:proofErr w:type="spellStart"/>'
$spellEnd = '<w:proofErr w:type="spellEnd"/>'

$paragraphXml =
    '<w:p w:rsidR="008F4B40" w:rsidRDefault="00527D01" w:rsidP="00527D01">' +
    '<w:pPr><w:ind w:firstLine="708"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
    (Run "Você é " $true) +
    $spellStart + (Run "Trevor" $false) + $spellEnd +
    (Run " Smith, ex-policial do departamento de polícia de Miami. " $true) +
    (Run "O ano é 2066. Já se passaram dois anos desde o falecimento de sua esposa e seu afastamento do trabalho. " $true) +
    $spellStart + (Run "Katelyn" $false) + $spellEnd +
    (Run " Parker, ou " $true) +
    $spellStart + (Run "Kat" $false) + $spellEnd +
    (Run ", como costumava ser chamada, foi alvo da vingança dos cartéis que na época travavam uma guerra ininterrupta com " $true) +
    $spellStart + (Run "Trevor" $false) + $spellEnd +
    (Run ". Após ser consumido pela perda e lamentação, você se vê em busca da única coisa que te trará paz... " $true) +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

$packageXml =
    '<?xml version="1.0" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $paragraphXml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

# InsertXML replaces the full contents of the target range in one shot, so
# the whole paragraph (minus its trailing paragraph mark) is rebuilt with
# the exact run/proofErr structure described above.
[void]$r.InsertXML($packageXml)
